$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "69.597.51"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "  +0.07%  "
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.Value = "2.506.14"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = "  -0.12%  "
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.Value = "  -0.04%  "
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.Value = "'575.23"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "  -0.29%  "
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.Value = "'166.19"
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.Value = "  -0.50%  "
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.Value = "'0.513"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = "  -0.69%  "
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.Value = "2.502.79"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "  -0.26%  "
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.Value = "'0.160"
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.Value = "  +0.51%  "
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.Value = "'0.168"
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.Value = "  -0.33%  "
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.Value = "'0.357"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.Value = "  +4.47%  "
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.Value = "'4.91"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = "  +1.03%  "
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.Value = "2.962.06"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = "  -0.08%  "
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.Value = "69.430.68"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = "  -0.07%  "
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = "  +1.09%  "
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.Value = "'24.84"
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.Value = "  -0.35%  "
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.Value = "2.507.28"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.Value = "  -0.18%  "
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.Value = "'11.22"
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = "  -1.96%  "
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.Value = "'7.54"
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.Value = "  -2.90%  "
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.Value = "'348.82"
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.Value = "  -0.45%  "
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = "  -1.25%  "
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.Value = "'1.95"
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.Value = "  -1.23%  "
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.Value = "  -0.09%  "
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.Value = "'70.19"
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.Value = "  +1.68%  "
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.Value = "'3.96"
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.Value = "  -1.36%  "
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.Value = "'8.80"
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.Value = "  -2.23%  "
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.Value = "  +0.01%  "
$c.Style = "Normal"
$c = $ws.Range("D29")
$c.Value = "'0.996"
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.Value = "  -0.30%  "
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.Value = "0.0₃0888"
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.Value = "  -1.68%  "
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.Value = "'7.84"
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.Value = "  -0.79%  "
$c.Style = "Normal"
$c = $ws.Range("D32")
$c.Value = "'460.09"
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.Value = "  -3.76%  "
$c.Style = "Normal"
$c = $ws.Range("D33")
$c.Value = "'1.21"
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.Value = "  -1.16%  "
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.Value = "  -0.08%  "
$c.Style = "Normal"
$c = $ws.Range("D36")
$c.Value = "'160.13"
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.Value = "  +3.90%  "
$c.Style = "Normal"
$c = $ws.Range("D37")
$c.Value = "'0.116"
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.Value = "  +0.00%  "
$c.Style = "Normal"
$c = $ws.Range("D38")
$c.Value = "'19.05"
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.Value = "  +0.79%  "
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.Value = "'18.45"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = "  -0.61%  "
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = "  +0.01%  "
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.Value = "  -0.32%  "
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = "  -1.78%  "
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.Value = "'1.59"
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = "  -0.91%  "
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.Value = "'38.18"
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = "  +0.15%  "
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.Value = "  -5.25%  "
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.Value = "  -7.40%  "
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.Value = "'142.26"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.Value = "  -1.36%  "
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.Value = "'3.46"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = "  -2.35%  "
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.Value = "'0.518"
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.Value = "  -2.48%  "
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.Value = "  +0.54%  "
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.Value = "'0.577"
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.Value = "  -1.76%  "
$c.Style = "Normal"
